$d = $word.ActiveDocument
try {
$d.XMLSchemaReferences.Add("http://schemas.microsoft.com/office/2006/metadata/contentType")
Write-Output "added schema ref"
} catch { Write-Output "ERR: $_" }
Write-Output $d.XMLSchemaReferences.Count
